$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.777.15'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.59%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.947.44'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '553.65'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.31'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +10.18%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +4.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.945.20'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.27%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.93%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.20%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +4.27%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000220'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.76'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +4.43%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.98%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.441.44'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.97'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +7.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.946.51'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '57.762.23'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '416.58'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.42'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +5.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.702'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +7.98%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.37'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +6.73%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.53'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +3.85%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.49'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.53%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +7.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.51'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +6.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.50'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.93'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0963'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.82%  '
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'Mantle'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.957'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +7.43%  '
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.71'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +7.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.05'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0₃0700'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +14.34%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +7.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '48.15'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.69'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +15.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '384.76'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +6.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.107'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0347'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.63%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.715.17'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '124.89'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +6.07%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +4.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.98'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.43%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.80'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.32%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.18%  '
